$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prodfix")

$newText = "UtilityOutcome - PRODFix_QOL_ECON - 9/19/2022"

$ws.Range("B2").Value = $newText
$ws.Range("B5").Value = $newText
$ws.Range("B8").Value = $newText
$ws.Range("B11").Value = $newText

$ws.Columns.Item(2).ColumnWidth = 42.25

$ws.Range("B11").Select()
